$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Insert a new row at position 13 which shifts everything below down by one.
$ws.Rows("13:13").Insert()

# Fill in the new Step2 row (row 13)
$ws.Range("D13").Value = "Step2"
$ws.Range("E13").Value = "'= (String[]) null"

# Update formula text in result1 row (now row 15)
$ws.Range("E15").Value = "'=flatten(null).length"

# Update formula text in result row (now row 16)
$ws.Range("E16").Value = "'=flatten(`$Step1).length"

# Fill in the new result2 row (row 17, which was previously an empty gap row)
$ws.Range("D17").Value = "result2"
$ws.Range("E17").Value = "'= flatten(`$Step2).length"

# Add the new result2 column (G) to the test table header rows (21 and 22)
$ws.Range("G21").Value = "_res_.`$result2"
$ws.Range("G22").Value = "_res_.`$result2"

# Add the new test data row (23)
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0

# Match the final cursor/selection position recorded in the workbook.
$null = $ws.Range("M14").Select()
